# vault backup: 2025-12-19 09:00:33
# Fills in the "第七天" (Day 7) column (H) of the sleep-diary table for
# all four weekly blocks' final week-block (rows 63-76), matching the
# values that were recorded for 2025-12-19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 63: "您今天早上几点醒来?" (wake-up time) - numeric time-of-day value,
# formatted the same way as the other day columns (B63/C63/... use h:mm).
$ws.Range("H63").Value = 0.28611111111111109
$ws.Range("H63").NumberFormat = $ws.Range("B63").NumberFormat

# Row 64: "您今天几点起床?" (get-up time) - recorded as text, like the
# other day columns in this table.
$ws.Range("H64").Value = "6：52"

# Row 65: "您昨晚几点上床?" (time went to bed)
$ws.Range("H65").Value = "23：05"

# Row 66: "您昨晚几点熄灯?" (lights-off time)
$ws.Range("H66").Value = "23：05"

# Row 67: "您昨晚熄灯后花了多长时间入睡(分钟)?" (minutes to fall asleep)
$ws.Range("H67").Value = 30

# Row 68: "您整晚醒来几次?" (number of times woke up)
$ws.Range("H68").Value = 1

# Row 69: "您整晚总共醒了多长时间(分钟)?" (total minutes awake)
$ws.Range("H69").Value = 5

# Row 70: "您整晚总共睡了多长时间(分钟)?" (total minutes slept)
$ws.Range("H70").Value = 440

# Row 71: substance use before bed - "无" (none)
$ws.Range("H71").Value = "无"

# Row 72: electronics use before bed - minutes
$ws.Range("H72").Value = 30

# Row 73: sleep quality rating
$ws.Range("H73").Value = 4

# Row 74: physical tension rating
$ws.Range("H74").Value = 3

# Row 75: mental tension rating
$ws.Range("H75").Value = 4

# Row 76: nap duration yesterday - "无" (none)
$ws.Range("H76").Value = "无"

# Match the workbook's final on-screen selection (the last cell the
# author edited).
$ws.Range("H76").Select()
